$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 551, shifting existing rows 551-671 down to 552-672
$ws.Rows.Item(551).Insert()

# Populate the newly inserted row 551 with the new weekly data record
$ws.Range("A551").Value = 6
$ws.Range("B551").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C551").Value = "Metropolitana"
$ws.Range("D551").Value = 44641
$ws.Range("E551").Value = 13
$ws.Range("F551").Value = 100112008
$ws.Range("G551").Value = "Coliflor"
$ws.Range("H551").Value = "Sin especificar"
$ws.Range("I551").Value = "Segunda"
$ws.Range("J551").Value = 2900
$ws.Range("K551").Value = 900
$ws.Range("L551").Value = 900
$ws.Range("M551").Value = 900
$ws.Range("N551").Value = "$/unidad"
$ws.Range("O551").Value = "Región Metropolitana"
$ws.Range("P551").Value = 900
$ws.Range("Q551").Value = 1
$ws.Range("R551").Value = "Hortaliza"
